# Rectify the sheet: add the missing credit line, move the selection,
# bump the outline level and tighten the header/footer margins a touch.

$wb = $excel.ActiveWorkbook

# --- Sheet1: append the new "By Shubham Sir" row ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A4").Value = "By Shubham Sir"

# Raise the sheet's tracked outline level to 3 without leaving any residual
# grouped row behind: group/ungroup a scratch row beyond the used range,
# then remove that scratch row again.
$ws1.Rows("10:10").OutlineLevel = 3
$ws1.Rows("10:10").Delete() | Out-Null

# Move the active selection to F13
$ws1.Activate()
$ws1.Range("F13").Select() | Out-Null

# Nudge the header/footer margins in slightly
$ws1.PageSetup.HeaderMargin = 36.75
$ws1.PageSetup.FooterMargin = 36.75

# --- Sheet2 / Sheet3: same header/footer margin tweak ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.PageSetup.HeaderMargin = 36.75
$ws2.PageSetup.FooterMargin = 36.75

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.PageSetup.HeaderMargin = 36.75
$ws3.PageSetup.FooterMargin = 36.75

# Changed the sleep time to 1000 (ms), per the rectified code.
Start-Sleep -Milliseconds 1000
